# Update "想去人数" (want-to-go count) values in column F
# for the "展览" (Exhibition) sheet and the "全部类型" (All Types) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 61
$ws1.Range("F4").Value  = 3601
$ws1.Range("F5").Value  = 2220
$ws1.Range("F7").Value  = 3
$ws1.Range("F8").Value  = 174
$ws1.Range("F9").Value  = 81
$ws1.Range("F10").Value = 69
$ws1.Range("F11").Value = 1329
$ws1.Range("F12").Value = 238
$ws1.Range("F13").Value = 1929
$ws1.Range("F14").Value = 139

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 61
$ws4.Range("F4").Value  = 3601
$ws4.Range("F5").Value  = 2220
$ws4.Range("F7").Value  = 3
$ws4.Range("F9").Value  = 174
$ws4.Range("F10").Value = 81
$ws4.Range("F11").Value = 69
$ws4.Range("F14").Value = 1329
$ws4.Range("F15").Value = 238
$ws4.Range("F16").Value = 1929
$ws4.Range("F17").Value = 139

$wb.Save()
